# LeadAssesment.xlsx - "Add files via upload" edit
# Updates TestScenario_1 / TestCase_1 (Step 1 and Step 2) data and marks
# the matching row in TestScenario_2 as a "Modified Step" as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2 (Step 1 : "Url") ---
# ExpectedResult text tweak
$ws.Range("H2").Value = "Url should be launched True"
# Flip Approved -> Rejected and add a reason
$ws.Range("I2").Value = "Rejected"
$ws.Range("J2").Value = "testingg"

# --- Row 3 (Step 2 : "Enter UserName") ---
$ws.Range("G3").Value = "Enter UserNames232"
$ws.Range("H3").Value = "text should be ablee to enter User Name "
$ws.Range("K3").Value = "Modified Step"

# --- Row 20 (TestScenario_2, Step 2 : "Enter UserName") ---
$ws.Range("K20").Value = "Modified Step"

# --- Update the active sheet view / selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("I10").Select()
